$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row data: row, open_price(D), close_price(E), high_price(F), low_price(G), shares_outstanding(H)
# fixed_ticker (I) is set to "OSS" for every data row (corrects a ticker-mapping bug from extra files).
$rows = @(
    @(2, 5.96999979019165, 4.783999919891357, 6, 4.449999809265137, 24478838),
    @(3, 5.96999979019165, 4.783999919891357, 6, 4.449999809265137, 24478838),
    @(4, 5.96999979019165, 4.783999919891357, 6, 4.449999809265137, 24478838),
    @(5, 5.96999979019165, 4.783999919891357, 6, 4.449999809265137, 24478838),
    @(6, 5.96999979019165, 4.783999919891357, 6, 4.449999809265137, 24478838),
    @(7, 5.96999979019165, 4.783999919891357, 6, 4.449999809265137, 24478838),
    @(8, 4.389999866485596, 4.46999979019165, 4.949999809265137, 4.099999904632568, 24478838),
    @(9, 4.150000095367432, 4.150000095367432, 4.75, 3.170000076293945, 24478838),
    @(10, 3.700000047683716, 2.470000028610229, 3.990000009536743, 2.25, 24478838),
    @(11, 1.929999947547913, 2.480000019073486, 2.700000047683716, 1.870000004768372, 24478838),
    @(12, 2.089999914169312, 2.329999923706055, 2.75, 1.950000047683716, 24478838),
    @(13, 1.730000019073486, 1.659999966621399, 2.450000047683716, 1.350000023841858, 24478838),
    @(14, 2.957000017166138, 2.299999952316284, 3.25, 2.25, 24478838),
    @(15, 2.039999961853028, 2.730000019073486, 3.009999990463257, 1.929999947547913, 24478838),
    @(16, 1.450000047683716, 1.590000033378601, 1.879999995231628, 1.149999976158142, 24478838),
    @(17, 2.009999990463257, 2.160000085830688, 2.940000057220459, 1.807999968528748, 24478838),
    @(18, 2.079999923706055, 2.240000009536743, 2.484999895095825, 2.079999923706055, 24478838),
    @(19, 4.099999904632568, 3.75, 4.980000019073486, 3.609999895095825, 24478838),
    @(20, 6.619999885559082, 5.320000171661377, 6.940000057220459, 4.409999847412109, 24478838),
    @(21, 5.800000190734863, 5.989999771118164, 6.189000129699707, 4.775000095367432, 24478838),
    @(22, 5.099999904632568, 5.010000228881836, 5.239999771118164, 4.699999809265137, 24478838),
    @(23, 4.909999847412109, 4.429999828338623, 5.190000057220459, 3.744999885559082, 24478838),
    @(24, 3.869999885559082, 4.699999809265137, 4.949999809265137, 3.799999952316284, 24478838),
    @(25, 3.950000047683716, 3.940000057220459, 3.980000019073486, 3.849999904632568, 24478838),
    @(26, 3.200000047683716, 3.220000028610229, 3.400000095367432, 2.960000038146973, 24478838),
    @(27, 3.009999990463257, 3.529999971389771, 3.710000038146973, 2.910000085830688, 24478838),
    @(28, 2.509999990463257, 2.470000028610229, 2.700000047683716, 2.289999961853028, 24478838),
    @(29, 2.900000095367432, 2.839999914169312, 3, 2.75, 24478838),
    @(30, 1.860000014305115, 2.029999971389771, 2.190000057220459, 1.800999999046326, 24478838),
    @(31, 2.140000104904175, 2.089999914169312, 2.25, 1.789999961853027, 24478838),
    @(32, 3.25, 2.839999914169312, 3.5, 2.545000076293945, 24478838),
    @(33, 2.400000095367432, 2.220000028610229, 2.559999942779541, 2.160000085830688, 24478838),
    @(34, 2.380000114440918, 2.299999952316284, 2.880000114440918, 2.220000028610229, 24478838),
    @(35, 3.410000085830688, 3.660000085830689, 4.329999923706055, 3.049999952316284, 24478838),
    @(36, 2.299999952316284, 2.309999942779541, 2.466000080108643, 1.850000023841858, 24478838),
    @(37, 3.549999952316284, 5.389999866485596, 6.400000095367432, 3.528000116348267, 24478838)
)

foreach ($row in $rows) {
    $r = $row[0]
    $ws.Range("D$r").Value = $row[1]
    $ws.Range("E$r").Value = $row[2]
    $ws.Range("F$r").Value = $row[3]
    $ws.Range("G$r").Value = $row[4]
    $ws.Range("H$r").Value = $row[5]
    $ws.Range("I$r").Value = "OSS"
}
